$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - rows 2-8, column F ("想去人数" = "want to go" count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 46
$ws1.Range("F3").Value = 128
$ws1.Range("F4").Value = 169
$ws1.Range("F5").Value = 3269
$ws1.Range("F6").Value = 330
$ws1.Range("F7").Value = 12
$ws1.Range("F8").Value = 416

# Sheet "全部类型" (all types) - same metric, rows offset by the extra "演出" rows 7-8
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 46
$ws4.Range("F3").Value = 128
$ws4.Range("F4").Value = 169
$ws4.Range("F5").Value = 3269
$ws4.Range("F6").Value = 330
$ws4.Range("F9").Value = 12
$ws4.Range("F10").Value = 416
